# Atualizacao de bases das ligas, do dia: 07-04-2024 as 22:30
# Hungary NB I: insert 2 new fixtures and refresh odds for the rows that shift down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Set-RowValues($r, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $val = $vals[$i]
        if ($val -ne $null) {
            $ws.Range($cols[$i] + $r).Value = $val
        }
    }
}

# Insert two blank rows right before the current row 157 so the three
# existing fixtures (old rows 157-159) shift down to rows 159-161.
$ws.Rows.Item(157).EntireRow.Insert()
$ws.Rows.Item(157).EntireRow.Insert()

# The insert borrows formatting from neighbouring rows but drops the
# bold/border style on column A for the freshly inserted rows - restore it
# explicitly so it matches the rest of the id column.
foreach ($r in 157, 158) {
    $cell = $ws.Range("A" + $r)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
}

# Row 157: brand-new fixture, MOL Fehervar FC vs Mezokovesd Zsory
Set-RowValues 157 @(155, 6818356, "Hungary NB I", "Hungary NB I", 45388.38541666666, "MOL Fehervar FC", "Mezokovesd Zsory", 5, 0, "H", 1.571, 3.6, 5, 1.666, 3.5, 4.75, -0.75, 1.875, 1.975, 2.5, 1.825, 2.025, 0.6659999999999999, -1, -1, 0.875, -1, 0.825, -1)

# Row 158: brand-new fixture, Diosgyori VTK vs Kisvarda FC
Set-RowValues 158 @(156, 6818354, "Hungary NB I", "Hungary NB I", 45388.47916666666, "Diosgyori VTK", "Kisvarda FC", 1, 1, "D", 1.727, 3.4, 4.2, 1.909, 3.3, 3.75, -0.5, 2, 1.85, 2.75, 1.975, 1.875, -1, 2.3, -1, -1, 0.8500000000000001, -1, 0.875)

# Row 159 (previously row 157): Debreceni VSC vs MTK Budapest, now with
# final-result columns H/I/J filled in and live odds instead of zeros.
Set-RowValues 159 @(157, 6818355, "Hungary NB I", "Hungary NB I", 45388.60416666666, "Debreceni VSC", "MTK Budapest", 1, 2, "A", 2, 3.4, 3.1, 1.666, 3.75, 4.5, -0.75, 1.925, 1.925, 2.75, 1.925, 1.925, -1, -1, 3.5, -1, 0.925, 0.4625, -0.5)

# Row 160 (previously row 158): Ferencvarosi TC vs Paksi, refreshed closing
# odds plus the final-result columns.
Set-RowValues 160 @(158, 6818352, "Hungary NB I", "Hungary NB I", 45389.40625, "Ferencvarosi TC", "Paksi", 1, 0, "H", 1.444, 4, 5.75, 1.4, 4.333, 7, -1.25, 1.9, 1.95, 3, 1.875, 1.975, 0.3999999999999999, -1, -1, -0.5, 0.475, -1, 0.9750000000000001)

# Row 161 (previously row 159): Kecskemeti TE vs Puskas Academy, refreshed
# closing odds plus the final-result columns.
Set-RowValues 161 @(159, 6818357, "Hungary NB I", "Hungary NB I", 45389.52083333334, "Kecskemeti TE", "Puskas Academy", 1, 2, "A", 2.4, 3.2, 2.6, 3.75, 3.5, 1.833, 0.5, 2.025, 1.825, 2.5, 2, 1.85, -1, -1, 0.833, -1, 0.825, 1, -1)
